$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.03132899571210146
$ws.Range("E2").Value = 0.8465258019277826
$ws.Range("G2").Value = 0.06328552949707955
$ws.Range("H2").Value = 0.5563925282331184
$ws.Range("I2").Value = 0.0363466510316357
$ws.Range("J2").Value = 0.09707645187154412
$ws.Range("K2").Value = 0.02347082155756652
$ws.Range("C3").Value = 715
$ws.Range("D3").Value = 0.0007684332085773349
$ws.Range("E3").Value = 0.2694176440127194
$ws.Range("F3").Value = 715
$ws.Range("G3").Value = 0.02111801982391626
$ws.Range("H3").Value = 0.1761669169645756
$ws.Range("I3").Value = 0.006703459774143994
$ws.Range("J3").Value = 0.03396003588568419
$ws.Range("K3").Value = 0.007764974609017372
$ws.Range("C4").Value = 715
$ws.Range("D4").Value = 0.01371360733173788
$ws.Range("E4").Value = 0.3480011560022831
$ws.Range("F4").Value = 715
$ws.Range("G4").Value = 0.02644675597548485
$ws.Range("H4").Value = 0.2291291250148788
$ws.Range("I4").Value = 0.009758537402376533
$ws.Range("J4").Value = 0.04374468501191586
$ws.Range("K4").Value = 0.009559829370118678
$ws.Range("C5").Value = 715
$ws.Range("D5").Value = 0.001299587893299758
$ws.Range("E5").Value = 0.3111274240072817
$ws.Range("F5").Value = 715
$ws.Range("G5").Value = 0.02459960756823421
$ws.Range("H5").Value = 0.2021545793395489
$ws.Range("I5").Value = 0.009164535324089229
$ws.Range("J5").Value = 0.03857860108837485
$ws.Range("K5").Value = 0.009044826379977167
$ws.Range("E6").Value = 1.07732873596251
$ws.Range("D7").Value = 0.02745914075057954
$ws.Range("E7").Value = 0.7452965029515326
$ws.Range("G7").Value = 0.05512112902943045
$ws.Range("H7").Value = 0.488878053962253
$ws.Range("I7").Value = 0.03434200666379184
$ws.Range("J7").Value = 0.0852981461212039
$ws.Range("K7").Value = 0.02053049392998219
$ws.Range("C8").Value = 518
$ws.Range("D8").Value = 0.0006273738108575344
$ws.Range("E8").Value = 0.2022015820257366
$ws.Range("F8").Value = 518
$ws.Range("G8").Value = 0.01577108656056225
$ws.Range("H8").Value = 0.1308854530798271
$ws.Range("I8").Value = 0.006195415859110653
$ws.Range("J8").Value = 0.02601894538383931
$ws.Range("K8").Value = 0.005786096910014749
$ws.Range("C9").Value = 518
$ws.Range("D9").Value = 0.01210114860441536
$ws.Range("E9").Value = 0.2791609010891989
$ws.Range("F9").Value = 518
$ws.Range("G9").Value = 0.02025187818799168
$ws.Range("H9").Value = 0.1845295174280182
$ws.Range("I9").Value = 0.009838219033554196
$ws.Range("J9").Value = 0.03408847749233246
$ws.Range("K9").Value = 0.007446028175763786
$ws.Range("C10").Value = 518
$ws.Range("D10").Value = 0.0009174455190077424
$ws.Range("E10").Value = 0.2039270030800253
$ws.Range("F10").Value = 518
$ws.Range("G10").Value = 0.01543852221220732
$ws.Range("H10").Value = 0.1317497495329008
$ws.Range("I10").Value = 0.007365834666416049
$ws.Range("J10").Value = 0.02613346138969064
$ws.Range("K10").Value = 0.005712611950002611
$ws.Range("E11").Value = 0.883647444890812
$ws.Range("D12").Value = 0.03061900229658931
$ws.Range("E12").Value = 0.8513227259973064
$ws.Range("G12").Value = 0.06384283024817705
$ws.Range("H12").Value = 0.5600891294889152
$ws.Range("I12").Value = 0.03298931068275124
$ws.Range("J12").Value = 0.09828316757921129
$ws.Range("K12").Value = 0.02422631497029215
$ws.Range("C13").Value = 730
$ws.Range("D13").Value = 0.0007934162858873606
$ws.Range("E13").Value = 0.2836754629388452
$ws.Range("F13").Value = 730
$ws.Range("G13").Value = 0.0223809196613729
$ws.Range("H13").Value = 0.1851774787064642
$ws.Range("I13").Value = 0.007017820025794208
$ws.Range("J13").Value = 0.03524050302803516
$ws.Range("K13").Value = 0.008320488268509507
$ws.Range("C14").Value = 730
$ws.Range("D14").Value = 0.01028321927879006
$ws.Range("E14").Value = 0.3006821880117059
$ws.Range("F14").Value = 730
$ws.Range("G14").Value = 0.02312570333015174
$ws.Range("H14").Value = 0.1983913854928687
$ws.Range("I14").Value = 0.008141801808960736
$ws.Range("J14").Value = 0.03662345535121858
$ws.Range("K14").Value = 0.008495513233356178
$ws.Range("C15").Value = 730
$ws.Range("D15").Value = 0.001296792179346085
$ws.Range("E15").Value = 0.3124301059870049
$ws.Range("F15").Value = 730
$ws.Range("G15").Value = 0.02493431628681719
$ws.Range("H15").Value = 0.2024387046694756
$ws.Range("I15").Value = 0.009120996575802565
$ws.Range("J15").Value = 0.03920478466898203
$ws.Range("K15").Value = 0.009039454045705497
$ws.Range("E16").Value = 1.009621520992368
$ws.Range("D17").Value = 0.01865265902597457
$ws.Range("E17").Value = 0.5483171029482037
$ws.Range("G17").Value = 0.04096021328587085
$ws.Range("H17").Value = 0.3597754908259958
$ws.Range("I17").Value = 0.02196197852026671
$ws.Range("J17").Value = 0.06441310828085989
$ws.Range("K17").Value = 0.01543511194176972
$ws.Range("C18").Value = 455
$ws.Range("D18").Value = 0.0005149961216375232
$ws.Range("E18").Value = 0.1739185689948499
$ws.Range("F18").Value = 455
$ws.Range("G18").Value = 0.01389433315489441
$ws.Range("H18").Value = 0.1133973768446594
$ws.Range("I18").Value = 0.004485916229896247
$ws.Range("J18").Value = 0.02178795787040144
$ws.Range("K18").Value = 0.005020585376769304
$ws.Range("C19").Value = 455
$ws.Range("D19").Value = 0.008697166456840932
$ws.Range("E19").Value = 0.2488408359931782
$ws.Range("F19").Value = 455
$ws.Range("G19").Value = 0.0189931420609355
$ws.Range("H19").Value = 0.1638728112448007
$ws.Range("I19").Value = 0.006947503890842199
$ws.Range("J19").Value = 0.03112816507928073
$ws.Range("K19").Value = 0.006911113508976996
$ws.Range("C20").Value = 455
$ws.Range("D20").Value = 0.0007867485983297229
$ws.Range("E20").Value = 0.1742072280030698
$ws.Range("F20").Value = 455
$ws.Range("G20").Value = 0.01360937999561429
$ws.Range("H20").Value = 0.1127434713998809
$ws.Range("I20").Value = 0.005213191383518279
$ws.Range("J20").Value = 0.02255417127162218
$ws.Range("K20").Value = 0.004959819838404655
$ws.Range("E21").Value = 0.9897040039068088
$ws.Range("C22").Value = 2169
$ws.Range("D22").Value = 0.03284667024854571
$ws.Range("E22").Value = 0.8966879429062828
$ws.Range("F22").Value = 2169
$ws.Range("G22").Value = 0.06800082640256733
$ws.Range("H22").Value = 0.5866953654913232
$ws.Range("I22").Value = 0.03994058596435934
$ws.Range("J22").Value = 0.1019669460365549
$ws.Range("K22").Value = 0.02500429155770689
$ws.Range("C23").Value = 844
$ws.Range("D23").Value = 0.000826406991109252
$ws.Range("E23").Value = 0.2956173539860174
$ws.Range("F23").Value = 844
$ws.Range("G23").Value = 0.02358393406029791
$ws.Range("H23").Value = 0.1913351627299562
$ws.Range("I23").Value = 0.007521819206885993
$ws.Range("J23").Value = 0.03799583250656724
$ws.Range("K23").Value = 0.008596397819928825
$ws.Range("C24").Value = 844
$ws.Range("D24").Value = 0.01128954836167395
$ws.Range("E24").Value = 0.3228164649335667
$ws.Range("F24").Value = 844
$ws.Range("G24").Value = 0.02452081348747015
$ws.Range("H24").Value = 0.2126501726452261
$ws.Range("I24").Value = 0.009120725793763995
$ws.Range("J24").Value = 0.03968897426966578
$ws.Range("K24").Value = 0.009195050457492471
$ws.Range("C25").Value = 844
$ws.Range("D25").Value = 0.001329442718997598
$ws.Range("E25").Value = 0.3160389759577811
$ws.Range("F25").Value = 844
$ws.Range("G25").Value = 0.02509266068227589
$ws.Range("H25").Value = 0.2045535693177953
$ws.Range("I25").Value = 0.009361646138131618
$ws.Range("J25").Value = 0.0400876043131575
$ws.Range("K25").Value = 0.009019757620990276
$ws.Range("E26").Value = 0.870586880017072
$ws.Range("C27").Value = 1728
$ws.Range("D27").Value = 0.02818129360675812
$ws.Range("E27").Value = 0.7776300153462217
$ws.Range("F27").Value = 1728
$ws.Range("G27").Value = 0.05824210569262504
$ws.Range("H27").Value = 0.5103661136003211
$ws.Range("I27").Value = 0.03311610657256096
$ws.Range("J27").Value = 0.08940756397787482
$ws.Range("K27").Value = 0.0217334067914635
$ws.Range("C28").Value = 652.4
$ws.Range("D28").Value = 0.000706125283613801
$ws.Range("E28").Value = 0.2449661223916337
$ws.Range("F28").Value = 652.4
$ws.Range("G28").Value = 0.01934965865220874
$ws.Range("H28").Value = 0.1593924776650965
$ws.Range("I28").Value = 0.006384886219166219
$ws.Range("J28").Value = 0.03100065493490547
$ws.Range("K28").Value = 0.007097708596847952
$ws.Range("C29").Value = 652.4
$ws.Range("D29").Value = 0.01121693800669163
$ws.Range("E29").Value = 0.2999003092059865
$ws.Range("F29").Value = 652.4
$ws.Range("G29").Value = 0.02266765860840678
$ws.Range("H29").Value = 0.1977146023651585
$ws.Range("I29").Value = 0.008761357585899531
$ws.Range("J29").Value = 0.03705475144088268
$ws.Range("K29").Value = 0.008321506949141622
$ws.Range("C30").Value = 652.4
$ws.Range("D30").Value = 0.001126003381796181
$ws.Range("E30").Value = 0.2635461474070325
$ws.Range("F30").Value = 652.4
$ws.Range("G30").Value = 0.02073489734902978
$ws.Range("H30").Value = 0.1707280148519203
$ws.Range("I30").Value = 0.008045240817591549
$ws.Range("J30").Value = 0.03331172454636544
$ws.Range("K30").Value = 0.007555293967016041
$ws.Range("E31").Value = 0.9661777171539143